$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the crypto price/volume refresh diff.
# D-column values that parse as plain numbers need an explicit
# text NumberFormat first, otherwise Excel silently coerces them
# to numeric (dropping significant trailing/leading zeros), which
# would not match the original inline-string cell content.

$ws.Range("D2").Value = '67.970.70'
$ws.Range("E2").Value = '  +1.75%  '
$ws.Range("D3").Value = '3.335.01'
$ws.Range("E3").Value = '  +1.91%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("D9").Value = '3.331.45'
$ws.Range("E10").Value = '  +4.71%  '
$ws.Range("E11").Value = '  +1.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.79'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '688.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").Value = '3.879.36'
$ws.Range("E15").Value = '  +2.20%  '
$ws.Range("E16").Value = '  +1.86%  '
$ws.Range("D17").Value = '67.986.73'
$ws.Range("E17").Value = '  +1.69%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.118'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.335.31'
$ws.Range("E19").Value = '  +1.71%  '
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("E21").Value = '  +2.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.898'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.08%  '
$ws.Range("E25").Value = '  +0.87%  '
$ws.Range("E26").Value = '  +0.74%  '
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.03'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.57'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.89%  '
$ws.Range("E31").Value = '  +5.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '573.49'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.07%  '
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("D35").Value = '3.713.88'
$ws.Range("E35").Value = '  -4.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '57.04'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.29'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.132'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.77%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.56%  '
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.38%  '
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0677'
$ws.Range("E44").Value = '  +0.94%  '
$ws.Range("E45").Value = '  +2.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0407'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("E47").Value = '  +5.70%  '
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("E50").Value = '  -2.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '129.61'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.05%  '
